# Scheduled-runner style update: refresh market-price-derived columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ, i.e. columns
# H, I, J, K, L, M, N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 33 (item 5512)
$ws.Range("H33").Value = 451.29166
$ws.Range("I33").Value = 434.6154
$ws.Range("K33").Value = 434.6154
$ws.Range("M33").Value = -205.6154

# ALC row 87 (item 10651)
$ws.Range("H87").Value = 37848.555
$ws.Range("J87").Value = 37848.555
$ws.Range("L87").Value = 37848.555
$ws.Range("N87").Value = -40344.555

# ALC row 90 (item 10651)
$ws.Range("H90").Value = 37848.555
$ws.Range("J90").Value = 37848.555
$ws.Range("L90").Value = 113545.665
$ws.Range("N90").Value = -126025.665

# ALC row 129 (item 36115)
$ws.Range("H129").Value = 841
$ws.Range("J129").Value = 1001.4286
$ws.Range("L129").Value = 3004.2858
$ws.Range("N129").Value = -13004.2858

# ALC row 137 (item 44013)
$ws.Range("H137").Value = 4181.485
$ws.Range("I137").Value = 2244.7778
$ws.Range("J137").Value = 6505.533
$ws.Range("K137").Value = 6734.3334
$ws.Range("L137").Value = 19516.599
$ws.Range("M137").Value = -4184.3334
$ws.Range("N137").Value = -24616.599

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32 (item 44147)
$ws.Range("H32").Value = 5804.413
$ws.Range("I32").Value = 5162.8955
$ws.Range("J32").Value = 14999.5
$ws.Range("K32").Value = 5162.8955
$ws.Range("L32").Value = 14999.5
$ws.Range("M32").Value = -4875.8955
$ws.Range("N32").Value = -15573.5

# ARM row 45 (item 27714)
$ws.Range("H45").Value = 1592.2258
$ws.Range("I45").Value = 1590.0834
$ws.Range("J45").Value = 1599.5714
$ws.Range("K45").Value = 1590.0834
$ws.Range("L45").Value = 1599.5714
$ws.Range("M45").Value = -1213.0834
$ws.Range("N45").Value = -2353.5714

# ARM row 61 (item 43999)
$ws.Range("H61").Value = 8639.489
$ws.Range("I61").Value = 4505.8076
$ws.Range("J61").Value = 14296.105
$ws.Range("K61").Value = 4505.8076
$ws.Range("L61").Value = 14296.105
$ws.Range("M61").Value = -4293.8076
$ws.Range("N61").Value = -14720.105

# ARM row 122 (item 36168)
$ws.Range("H122").Value = 15628800
$ws.Range("I122").Value = 4116.8335
$ws.Range("J122").Value = 62502850
$ws.Range("K122").Value = 12350.5005
$ws.Range("L122").Value = 187508550
$ws.Range("M122").Value = -9900.500499999998
$ws.Range("N122").Value = -187513450

# ARM row 132 (item 43997)
$ws.Range("H132").Value = 5883.8203
$ws.Range("I132").Value = 2276.4375
$ws.Range("J132").Value = 8393.305
$ws.Range("K132").Value = 6829.3125
$ws.Range("L132").Value = 25179.915
$ws.Range("M132").Value = -4299.3125
$ws.Range("N132").Value = -30239.915

# ARM row 136 (item 43999)
$ws.Range("H136").Value = 8639.489
$ws.Range("I136").Value = 4505.8076
$ws.Range("J136").Value = 14296.105
$ws.Range("K136").Value = 13517.4228
$ws.Range("L136").Value = 42888.315
$ws.Range("M136").Value = -10967.4228
$ws.Range("N136").Value = -47988.315

$ws = $wb.Worksheets.Item("BSM")
# BSM row 63 (item 10592)
$ws.Range("H63").Value = 32000
$ws.Range("J63").Value = 32000
$ws.Range("L63").Value = 32000
$ws.Range("N63").Value = -33372

# BSM row 66 (item 10592)
$ws.Range("H66").Value = 32000
$ws.Range("J66").Value = 32000
$ws.Range("L66").Value = 96000
$ws.Range("N66").Value = -102864

# BSM row 112 (item 25788)
$ws.Range("H112").Value = 41989.668
$ws.Range("J112").Value = 41989.668
$ws.Range("L112").Value = 41989.668
$ws.Range("N112").Value = -44943.668

$ws = $wb.Worksheets.Item("CRP")
# CRP row 31 (item 44023)
$ws.Range("H31").Value = 2093.2126
$ws.Range("I31").Value = 1474.3279
$ws.Range("J31").Value = 3237.2122
$ws.Range("K31").Value = 1474.3279
$ws.Range("L31").Value = 3237.2122
$ws.Range("M31").Value = -1179.3279
$ws.Range("N31").Value = -3827.2122

# CRP row 34 (item 44023)
$ws.Range("H34").Value = 2093.2126
$ws.Range("I34").Value = 1474.3279
$ws.Range("J34").Value = 3237.2122
$ws.Range("K34").Value = 1474.3279
$ws.Range("L34").Value = 3237.2122
$ws.Range("M34").Value = -1272.3279
$ws.Range("N34").Value = -3641.2122

# CRP row 48 (item 3870)
$ws.Range("H48").Value = 12384.333
$ws.Range("J48").Value = 12384.333
$ws.Range("L48").Value = 12384.333
$ws.Range("N48").Value = -13336.333

# CRP row 107 (item 27689)
$ws.Range("H107").Value = 1158.5454
$ws.Range("I107").Value = 1290.1428
$ws.Range("J107").Value = 928.25
$ws.Range("K107").Value = 1290.1428
$ws.Range("L107").Value = 928.25
$ws.Range("M107").Value = 629.8571999999999
$ws.Range("N107").Value = -4768.25

# CRP row 110 (item 25791)
$ws.Range("H110").Value = 24851
$ws.Range("J110").Value = 24851
$ws.Range("L110").Value = 24851
$ws.Range("N110").Value = -33031

$ws = $wb.Worksheets.Item("CUL")
# CUL row 3 (item 44094)
$ws.Range("H3").Value = 8236.684999999999
$ws.Range("I3").Value = 5654.0835
$ws.Range("J3").Value = 12664
$ws.Range("K3").Value = 16962.2505
$ws.Range("L3").Value = 37992
$ws.Range("M3").Value = -16850.2505
$ws.Range("N3").Value = -38216

# CUL row 5 (item 43974)
$ws.Range("H5").Value = 6177356.5
$ws.Range("I5").Value = 660.0952
$ws.Range("J5").Value = 27795792
$ws.Range("K5").Value = 1980.2856
$ws.Range("L5").Value = 83387376
$ws.Range("M5").Value = -1868.2856
$ws.Range("N5").Value = -83387600

# CUL row 18 (item 36056)
$ws.Range("H18").Value = 14286364
$ws.Range("I18").Value = 18182118
$ws.Range("J18").Value = 1933.3334
$ws.Range("K18").Value = 54546354
$ws.Range("L18").Value = 5800.0002
$ws.Range("M18").Value = -54546185
$ws.Range("N18").Value = -6138.0002

# CUL row 23 (item 4858)
$ws.Range("H23").Value = 810.95
$ws.Range("I23").Value = 2025.3334
$ws.Range("J23").Value = 290.5
$ws.Range("K23").Value = 6076.0002
$ws.Range("L23").Value = 871.5
$ws.Range("M23").Value = -5841.0002
$ws.Range("N23").Value = -1341.5

# CUL row 68 (item 12895)
$ws.Range("H68").Value = 8500.308000000001
$ws.Range("I68").Value = 584.2857
$ws.Range("J68").Value = 17735.666
$ws.Range("K68").Value = 1752.8571
$ws.Range("L68").Value = 53206.99800000001
$ws.Range("M68").Value = -941.8571000000002
$ws.Range("N68").Value = -54828.99800000001

# CUL row 69 (item 12850)
$ws.Range("H69").Value = 22728164
$ws.Range("J69").Value = 25000932
$ws.Range("L69").Value = 75002796
$ws.Range("N69").Value = -75004418

# CUL row 71 (item 12895)
$ws.Range("H71").Value = 8500.308000000001
$ws.Range("I71").Value = 584.2857
$ws.Range("J71").Value = 17735.666
$ws.Range("K71").Value = 5258.571300000001
$ws.Range("L71").Value = 159620.994
$ws.Range("M71").Value = -1202.571300000001
$ws.Range("N71").Value = -167732.994

# CUL row 72 (item 12850)
$ws.Range("H72").Value = 22728164
$ws.Range("J72").Value = 25000932
$ws.Range("L72").Value = 225008388
$ws.Range("N72").Value = -225016500

# CUL row 75 (item 12863)
$ws.Range("H75").Value = 1500
$ws.Range("J75").Value = 1500
$ws.Range("L75").Value = 4500
$ws.Range("N75").Value = -6496

# CUL row 78 (item 12863)
$ws.Range("H78").Value = 1500
$ws.Range("J78").Value = 1500
$ws.Range("L78").Value = 13500
$ws.Range("N78").Value = -23484

# CUL row 103 (item 19839)
$ws.Range("H103").Value = 2257
$ws.Range("I103").Value = 2333.3333
$ws.Range("K103").Value = 6999.999899999999
$ws.Range("M103").Value = -6120.999899999999

# CUL row 110 (item 27857)
$ws.Range("H110").Value = 3998.6924
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 3998.6924
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 11996.0772
$ws.Range("M110").ClearContents()
$ws.Range("N110").Value = -20176.0772

# CUL row 113 (item 27843)
$ws.Range("H113").Value = 734.19275
$ws.Range("J113").Value = 631.6
$ws.Range("L113").Value = 1894.8
$ws.Range("N113").Value = -6234.8

# CUL row 135 (item 43974)
$ws.Range("H135").Value = 6177356.5
$ws.Range("I135").Value = 660.0952
$ws.Range("J135").Value = 27795792
$ws.Range("K135").Value = 5940.8568
$ws.Range("L135").Value = 250162128
$ws.Range("M135").Value = -3405.8568
$ws.Range("N135").Value = -250167198

$ws = $wb.Worksheets.Item("GSM")
# GSM row 80 (item 12521)
$ws.Range("H80").Value = 27000
$ws.Range("I80").Value = 27000
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 27000
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -26002
$ws.Range("N80").ClearContents()

# GSM row 83 (item 12521)
$ws.Range("H83").Value = 27000
$ws.Range("I83").Value = 27000
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 135000
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -130008
$ws.Range("N83").ClearContents()

# GSM row 132 (item 44008)
$ws.Range("H132").Value = 58740.65
$ws.Range("I132").Value = 254294.75
$ws.Range("J132").Value = 9852.125
$ws.Range("K132").Value = 762884.25
$ws.Range("L132").Value = 29556.375
$ws.Range("M132").Value = -760354.25
$ws.Range("N132").Value = -34616.375

$ws = $wb.Worksheets.Item("LTW")
# LTW row 118 (item 26146)
$ws.Range("H118").Value = 45000
$ws.Range("J118").Value = 45000
$ws.Range("L118").Value = 45000
$ws.Range("N118").Value = -48314

$ws = $wb.Worksheets.Item("WVR")
# WVR row 132 (item 44029)
$ws.Range("H132").Value = 5006.65
$ws.Range("I132").Value = 4404.9165
$ws.Range("J132").Value = 5909.25
$ws.Range("K132").Value = 6909.999999999999
$ws.Range("L132").Value = 17727.75
$ws.Range("M132").Value = -10684.7495
$ws.Range("N132").Value = -22787.75
